# Handle mid-quarter storage changes:
# Append a new form-response row (row 3) for Bob Banana, mirroring the
# existing Diane Durian row (row 2), including the mailto: hyperlink on
# the email-address cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A3: Timestamp -------------------------------------------------------
# Match the direct formatting already used by A2 (date/time number format,
# Arial 11 black) so no redundant style gets introduced.
$ws.Range("A3").Value = 44180.7005092593
$ws.Range("A3").NumberFormat = 'm/d/yyyy\ h:mm:ss'
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Color = 0

# --- B3: Email address (with hyperlink) -----------------------------------
$null = $ws.Hyperlinks.Add($ws.Range("B3"), "mailto:bbanana@example.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "bbanana@example.com")
# Reset the auto-applied "Hyperlink" character formatting back to the
# plain/default look used elsewhere on this sheet.
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").Font.Color = 0
$ws.Range("B3").Font.Underline = -4142

# --- Remaining row-3 fields -------------------------------------------------
$ws.Range("C3").Value = "Bob"
$ws.Range("D3").Value = "Banana"
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = "BBBB"
$ws.Range("G3").Value = "no"
$ws.Range("H3").Value = "Yes"

# Final selection left on G3, matching the saved view state.
$null = $ws.Range("G3").Select()
